$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for the refreshed crypto symbol list.
# Values are written with a leading apostrophe to preserve them as literal text
# (matching the original "inlineStr" cell type) rather than being parsed as
# numbers/percentages by Excel.
$ws.Range("D2").Value = "'256.17"
$ws.Range("E2").Value = "'0.38%"
$ws.Range("D3").Value = "'26.94"
$ws.Range("E3").Value = "'-4.27%"
$ws.Range("D4").Value = "'4.718"
$ws.Range("E4").Value = "'-10.16%"
$ws.Range("E5").Value = "'1.46%"
$ws.Range("D6").Value = "'6.659"
$ws.Range("E6").Value = "'-0.71%"
$ws.Range("D7").Value = "'0.8678"
$ws.Range("E7").Value = "'-0.05%"
$ws.Range("D8").Value = "'0.9581"
$ws.Range("E8").Value = "'-6.94%"
$ws.Range("D9").Value = "'0.1403"
$ws.Range("E9").Value = "'-0.67%"
$ws.Range("D10").Value = "'0.03949"
$ws.Range("E10").Value = "'13.37%"
$ws.Range("D11").Value = "'0.07155"
$ws.Range("E11").Value = "'0.30%"
$ws.Range("D12").Value = "'0.03199"
$ws.Range("E12").Value = "'0.46%"
$ws.Range("D13").Value = "'0.09251"
$ws.Range("E13").Value = "'0.27%"
$ws.Range("D14").Value = "'0.001548"
$ws.Range("E14").Value = "'-0.25%"
$ws.Range("D15").Value = "'0.0006062"
$ws.Range("E15").Value = "'0.34%"
$ws.Range("E16").Value = "'4.18%"
$ws.Range("E17").Value = "'-0.44%"
$ws.Range("D18").Value = "'3.203"
$ws.Range("E18").Value = "'-0.86%"
$ws.Range("E19").Value = "'-0.18%"
$ws.Range("D20").Value = "'0.3135"
$ws.Range("E20").Value = "'-1.50%"
$ws.Range("E21").Value = "'-1.41%"
$ws.Range("D22").Value = "'3.812"
$ws.Range("E22").Value = "'7.90%"
$ws.Range("D23").Value = "'0.04220"
$ws.Range("E23").Value = "'1.71%"
$ws.Range("D25").Value = "'0.001219"
$ws.Range("E25").Value = "'-0.02%"
$ws.Range("D26").Value = "'0.004491"
$ws.Range("E26").Value = "'-7.98%"
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("E27").Value = "'0.02%"
$ws.Range("D28").Value = "'0.0001938"
$ws.Range("E28").Value = "'-0.04%"
$ws.Range("D40").Value = "'0.03820"
$ws.Range("E40").Value = "'0.25%"
$ws.Range("D41").Value = "'0.006133"
$ws.Range("E41").Value = "'60.91%"
$ws.Range("D42").Value = "'0.1099"
$ws.Range("E42").Value = "'-0.19%"
$ws.Range("D43").Value = "'0.002201"
$ws.Range("E43").Value = "'-6.47%"
$ws.Range("D44").Value = "'0.01059"
$ws.Range("E44").Value = "'9.28%"
$ws.Range("D45").Value = "'0.00005502"
$ws.Range("E45").Value = "'5.03%"
$ws.Range("E46").Value = "'0.01%"
$ws.Range("D47").Value = "'0.08856"
$ws.Range("E47").Value = "'-4.81%"
$ws.Range("D48").Value = "'0.002393"
$ws.Range("E48").Value = "'11.21%"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("E49").Value = "'0.01%"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("E50").Value = "'0.01%"

Write-Output "Updated crypto price/volume cells."
